$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Delete the old end_mdr_introduce_time row (row 6); everything below shifts up.
$ws.Rows.Item(6).Delete()

# Update row 5 label cell (A5) and description cell (E5) to new text/order.
$ws.Range("A5").Value = "mdr_introduce_time"
$ws.Range("E5").Value = "Calendar year that MDR-TB first begins to emerge"

# Update selection to match target state.
$ws.Range("A5").Select()
